# Figs_4_5_ELP_estimation_error.xlsx - "Add files via upload" re-upload
# Net effect (per the OOXML diff):
#   - Column A (rows 2-13) on sheet "Hoja1" has every value's sign flipped.
#   - Cell A2 loses its wrap-text formatting override (back to the default style).
#   - The sheet's saved cursor/selection moves from C17 to C22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the wrap-text formatting that singled out A2 so it reverts to the
# workbook's default (unstyled) cell format.
$ws.Range("A2").ClearFormats()

# Flip the sign of every value in A2:A13.
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = -1 * $cell.Value2
}

# Restore the sheet's recorded selection/active cell.
$ws.Range("C22").Select()
